$wb = $excel.ActiveWorkbook

# --- Repayment Schedule sheet: insert a new blank column before column N ---
$repayWs = $wb.Worksheets.Item("Repayment Schedule")
$repayWs.Columns("N").Insert()

# Give the newly inserted column (N) the same width as column M ("Paid"/In Advance column)
$repayWs.Columns("N").ColumnWidth = $repayWs.Columns("M").ColumnWidth

# Update selection on the Repayment Schedule sheet and make it the active sheet/tab
$repayWs.Activate()
$repayWs.Range("L21").Select()

# --- Transactions sheet: restore its own selection (it is no longer the active tab) ---
$transWs = $wb.Worksheets.Item("Transactions")
$transWs.Range("E3").Select()

# Re-activate Repayment Schedule so it remains the selected/visible tab on save
$repayWs.Activate()
